$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. "1.002", "28.251.91").
# Force the whole Price column to Text format before writing so Excel keeps the
# exact original text representation (including trailing zeros / multi-dot values)
# instead of silently converting to a floating point number.
$ws.Range("D2:D51").NumberFormat = "@"

$priceUpdates = @{
    'D2' = '28.251.91'
    'D3' = '1.813.51'
    'D4' = '1.002'
    'D5' = '327.87'
    'D7' = '0.4337'
    'D9' = '44.84'
    'D10' = '0.07676'
    'D12' = '1.001'
    'D13' = '22.01'
    'D14' = '6.290'
    'D15' = '7.508'
    'D16' = '1.830.68'
    'D17' = '93.78'
    'D19' = '0.06560'
    'D21' = '17.48'
    'D22' = '6.248'
    'D23' = '28.274.75'
    'D25' = '2.060'
    'D26' = '162.70'
    'D27' = '20.63'
    'D28' = '2.034.35'
    'D29' = '2.280'
    'D30' = '128.60'
    'D31' = '1.208'
    'D32' = '5.937'
    'D33' = '0.09165'
    'D34' = '3.473'
    'D36' = '0.02350'
    'D37' = '0.2172'
    'D38' = '5.190'
    'D39' = '0.6567'
    'D40' = '0.06194'
    'D42' = '8.109'
    'D43' = '1.433'
    'D45' = '13.80'
    'D46' = '0.6098'
    'D47' = '3.750'
    'D48' = '125.52'
    'D49' = '2.014'
    'D50' = '1.155'
    'D51' = '0.07004'
}
foreach ($cell in $priceUpdates.Keys) {
    $ws.Range($cell).Value = $priceUpdates[$cell]
}

# Drop the temporary Text number format again so the cells keep using the
# default (General) style, matching the original file layout.
$ws.Range("D2:D51").Style = "Normal"

# Column E holds the Volume(1h) percentage strings (with surrounding spaces);
# these are never numeric-looking so a plain Value assignment keeps them as text.
$volumeUpdates = @{
    'E2' = '  +2.94%  '
    'E3' = '  +3.95%  '
    'E4' = '  -0.17%  '
    'E5' = '  +1.62%  '
    'E6' = '  -0.03%  '
    'E7' = '  +2.77%  '
    'E8' = '  +2.27%  '
    'E9' = '  -1.21%  '
    'E10' = '  +3.65%  '
    'E11' = '  +2.81%  '
    'E12' = '  -0.06%  '
    'E13' = '  +2.98%  '
    'E14' = '  +3.35%  '
    'E15' = '  +4.69%  '
    'E16' = '  +5.11%  '
    'E17' = '  +7.46%  '
    'E18' = '  +1.67%  '
    'E19' = '  +6.07%  '
    'E20' = '  +0.01%  '
    'E22' = '  +2.64%  '
    'E23' = '  +2.89%  '
    'E24' = '  -0.25%  '
    'E25' = '  -11.36%  '
    'E26' = '  +7.19%  '
    'E27' = '  +1.08%  '
    'E28' = '  +4.88%  '
    'E29' = '  -3.23%  '
    'E30' = '  +2.09%  '
    'E31' = '  +0.38%  '
    'E32' = '  +4.79%  '
    'E33' = '  +0.44%  '
    'E34' = '  -5.83%  '
    'E35' = '  +2.98%  '
    'E36' = '  +2.95%  '
    'E37' = '  +2.15%  '
    'E38' = '  +2.07%  '
    'E39' = '  +2.83%  '
    'E40' = '  +2.00%  '
    'E41' = '  +0.01%  '
    'E42' = '  +2.98%  '
    'E43' = '  +0.78%  '
    'E44' = '  -0.03%  '
    'E45' = '  +1.01%  '
    'E46' = '  +4.19%  '
    'E47' = '  +0.84%  '
    'E48' = '  +0.60%  '
    'E49' = '  +3.44%  '
    'E50' = '  +3.09%  '
    'E51' = '  +2.20%  '
}
foreach ($cell in $volumeUpdates.Keys) {
    $ws.Range($cell).Value = $volumeUpdates[$cell]
}
